$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 611
$ws.Range("I92").Value = 587.17645
$ws.Range("J92").Value = 692
$ws.Range("K92").Value = 587.17645
$ws.Range("L92").Value = 692
$ws.Range("M92").Value = 660.82355
$ws.Range("N92").Value = -3188

$ws.Range("H106").Value = 1598.25
$ws.Range("I106").Value = 1186
$ws.Range("K106").Value = 1186
$ws.Range("M106").Value = -555

$ws.Range("H137").Value = 43424.25
$ws.Range("I137").Value = 1390
$ws.Range("J137").Value = 73448.71000000001
$ws.Range("K137").Value = 4170
$ws.Range("L137").Value = 220346.13
$ws.Range("M137").Value = -1620
$ws.Range("N137").Value = -225446.13

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1397.2106
$ws.Range("I2").Value = 848.7692
$ws.Range("J2").Value = 2585.5
$ws.Range("K2").Value = 848.7692
$ws.Range("L2").Value = 2585.5
$ws.Range("M2").Value = -735.7692
$ws.Range("N2").Value = -2811.5

$ws.Range("H74").Value = 651.875
$ws.Range("I74").Value = 330.26086
$ws.Range("J74").Value = 1473.7778
$ws.Range("K74").Value = 330.26086
$ws.Range("L74").Value = 1473.7778
$ws.Range("M74").Value = 543.73914
$ws.Range("N74").Value = -3221.7778

$ws.Range("H77").Value = 651.875
$ws.Range("I77").Value = 330.26086
$ws.Range("J77").Value = 1473.7778
$ws.Range("K77").Value = 1651.3043
$ws.Range("L77").Value = 7368.889
$ws.Range("M77").Value = 2716.6957
$ws.Range("N77").Value = -16104.889

$ws.Range("H116").Value = 1397.2106
$ws.Range("I116").Value = 848.7692
$ws.Range("J116").Value = 2585.5
$ws.Range("K116").Value = 848.7692
$ws.Range("L116").Value = 2585.5
$ws.Range("M116").Value = 1445.2308
$ws.Range("N116").Value = -7173.5

$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1397.2106
$ws.Range("I3").Value = 848.7692
$ws.Range("J3").Value = 2585.5
$ws.Range("K3").Value = 848.7692
$ws.Range("L3").Value = 2585.5
$ws.Range("M3").Value = -734.7692
$ws.Range("N3").Value = -2813.5

$ws.Range("H22").Value = 320.85715
$ws.Range("I22").Value = 320.85715
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 320.85715
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -147.85715
$ws.Range("N22").ClearContents()

$ws.Range("H107").Value = 1005.0714
$ws.Range("I107").Value = 1005.46155
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1005.46155
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 914.53845
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1242.8572
$ws.Range("J16").Value = 1300
$ws.Range("L16").Value = 1300
$ws.Range("N16").Value = -1874

$ws.Range("H31").Value = 10759.846
$ws.Range("I31").Value = 15082.218
$ws.Range("J31").Value = 4546.4375
$ws.Range("K31").Value = 15082.218
$ws.Range("L31").Value = 4546.4375
$ws.Range("M31").Value = -14787.218
$ws.Range("N31").Value = -5136.4375

$ws.Range("H34").Value = 10759.846
$ws.Range("I34").Value = 15082.218
$ws.Range("J34").Value = 4546.4375
$ws.Range("K34").Value = 15082.218
$ws.Range("L34").Value = 4546.4375
$ws.Range("M34").Value = -14880.218
$ws.Range("N34").Value = -4950.4375

$ws.Range("H99").Value = 5642.222
$ws.Range("I99").Value = 4056
$ws.Range("K99").Value = 4056
$ws.Range("M99").Value = -2558

$ws.Range("H113").Value = 1242.8572
$ws.Range("J113").Value = 1300
$ws.Range("L113").Value = 1300
$ws.Range("N113").Value = -5640

$ws.Range("H126").Value = 5642.222
$ws.Range("I126").Value = 4056
$ws.Range("K126").Value = 12168
$ws.Range("M126").Value = -9698

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 823.1
$ws.Range("J34").Value = 899
$ws.Range("L34").Value = 2697
$ws.Range("N34").Value = -2865

$ws.Range("H107").Value = 4152.0356
$ws.Range("I107").Value = 8175.769
$ws.Range("J107").Value = 664.8
$ws.Range("K107").Value = 24527.307
$ws.Range("L107").Value = 1994.4
$ws.Range("M107").Value = -22607.307
$ws.Range("N107").Value = -5834.4

$ws.Range("H129").Value = 278945.84
$ws.Range("I129").Value = 503.125
$ws.Range("J129").Value = 501700
$ws.Range("K129").Value = 1509.375
$ws.Range("L129").Value = 1505100
$ws.Range("M129").Value = 3490.625
$ws.Range("N129").Value = -1515100

$ws.Range("H131").Value = 179420.12
$ws.Range("J131").Value = 193143.6
$ws.Range("L131").Value = 579430.8
$ws.Range("N131").Value = -589510.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4789.8
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 4987.25
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 4987.25
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -5527.25

$ws.Range("H73").Value = 4789.8
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 4987.25
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 4987.25
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -6859.25

$ws.Range("H80").Value = 10785.667
$ws.Range("I80").Value = 22841
$ws.Range("J80").Value = 4758
$ws.Range("K80").Value = 22841
$ws.Range("L80").Value = 4758
$ws.Range("M80").Value = -21843
$ws.Range("N80").Value = -6754

$ws.Range("H83").Value = 10785.667
$ws.Range("I83").Value = 22841
$ws.Range("J83").Value = 4758
$ws.Range("K83").Value = 114205
$ws.Range("L83").Value = 23790
$ws.Range("M83").Value = -109213
$ws.Range("N83").Value = -33774

$ws.Range("H107").Value = 1561.5555
$ws.Range("I107").Value = 394
$ws.Range("J107").Value = 3896.6667
$ws.Range("K107").Value = 394
$ws.Range("L107").Value = 3896.6667
$ws.Range("M107").Value = 1526
$ws.Range("N107").Value = -7736.6667

$ws.Range("H132").Value = 43208.418
$ws.Range("I132").Value = 35625.863
$ws.Range("J132").Value = 74621.86
$ws.Range("K132").Value = 106877.589
$ws.Range("L132").Value = 223865.58
$ws.Range("M132").Value = -104347.589
$ws.Range("N132").Value = -228925.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1450
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -105
$ws.Range("N22").Value = -3090

$ws.Range("H27").Value = 1450
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -293
$ws.Range("N27").Value = -2714

$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("K30").Value = 500
$ws.Range("M30").Value = -392

$ws.Range("H132").Value = 525731.4
$ws.Range("I132").Value = 928182.1
$ws.Range("K132").Value = 2784546.3
$ws.Range("M132").Value = -2782016.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 257507
$ws.Range("J26").Value = 257507
$ws.Range("L26").Value = 257507
$ws.Range("N26").Value = -258093

$ws.Range("H81").Value = 1471.4286
$ws.Range("I81").Value = 1550
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 3100
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -2039
$ws.Range("N81").Value = -4122

$ws.Range("H84").Value = 1471.4286
$ws.Range("I84").Value = 1550
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 15500
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -10196
$ws.Range("N84").Value = -20608

$ws.Range("H122").Value = 2183.25
$ws.Range("I122").Value = 1966.6666
$ws.Range("K122").Value = 5899.9998
$ws.Range("M122").Value = -3449.9998

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 1217.32
$ws.Range("I132").Value = 962.6923
$ws.Range("J132").Value = 2120.0908
$ws.Range("K132").Value = 2888.0769
$ws.Range("L132").Value = 6360.2724
$ws.Range("M132").Value = -358.0769
$ws.Range("N132").Value = -11420.2724
